{"js": "const paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst firstParagraph = paragraphs.items[0];\nfirstParagraph.insertParagraph(\"mmyyyyyyyy\", Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$firstPara = $d.Paragraphs.Item(1)\n$firstPara.Range.InsertParagraphAfter()\n$newPara = $d.Paragraphs.Item(2)\n$newPara.Range.Text = \"mmyyyyyyyy\"\n"}
